$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 3: Fecha D3 44405 -> 44424
$ws.Range("D3").Value = 44424

# Update existing row 4: Fecha D4 44343 -> 44405, and Volumen/Precios change
$ws.Range("D4").Value = 44405
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 1200
$ws.Range("O4").Value = 1200
$ws.Range("P4").Value = 1200
$ws.Range("S4").Value = 1200

# New row 5
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "Vega Modelo de Temuco"
$ws.Range("C5").Value = "La Araucanía"
$ws.Range("D5").Value = 44417
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100108
$ws.Range("H5").Value = "Tropicales y subtropicales"
$ws.Range("I5").Value = 100108001
$ws.Range("J5").Value = "Guayaba"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 1200
$ws.Range("O5").Value = 1200
$ws.Range("P5").Value = 1200
$ws.Range("Q5").Value = "$/kilo"
$ws.Range("R5").Value = "Región de Arica y Parinacota"
$ws.Range("S5").Value = 1200
$ws.Range("T5").Value = 1

# New row 6 (carries the former row-4 data: D=44343, M=60, N/O/P=1300, S=1300)
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "Vega Modelo de Temuco"
$ws.Range("C6").Value = "La Araucanía"
$ws.Range("D6").Value = 44343
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100108
$ws.Range("H6").Value = "Tropicales y subtropicales"
$ws.Range("I6").Value = 100108001
$ws.Range("J6").Value = "Guayaba"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 1300
$ws.Range("O6").Value = 1300
$ws.Range("P6").Value = 1300
$ws.Range("Q6").Value = "$/kilo"
$ws.Range("R6").Value = "Región de Arica y Parinacota"
$ws.Range("S6").Value = 1300
$ws.Range("T6").Value = 1

# New row 7
$ws.Range("A7").Value = 10
$ws.Range("B7").Value = "Vega Modelo de Temuco"
$ws.Range("C7").Value = "La Araucanía"
$ws.Range("D7").Value = 44418
$ws.Range("E7").Value = 9
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100108
$ws.Range("H7").Value = "Tropicales y subtropicales"
$ws.Range("I7").Value = 100108001
$ws.Range("J7").Value = "Guayaba"
$ws.Range("K7").Value = "Sin especificar"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = 1200
$ws.Range("O7").Value = 1200
$ws.Range("P7").Value = 1200
$ws.Range("Q7").Value = "$/kilo"
$ws.Range("R7").Value = "Región de Arica y Parinacota"
$ws.Range("S7").Value = 1200
$ws.Range("T7").Value = 1

# Apply the date format (style index 2 in before.xlsx / numFmtId 165) to the new D cells
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
